$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the second row of teacher data (row 2) with new values.
# Order matches the original shared-string slot order (name, phone, email, maGV)
# so the rebuilt shared strings table lines up with the target layout.
$ws.Range("B2").Value = "GV Thịnh Nguyễn"
$ws.Range("D2").Value = "0383965078"
# E2 already carries the Hyperlink/quote-prefixed text style; prefix with an
# apostrophe (standard Excel "force text" marker) so the cell keeps that
# style instead of being reformatted when the value is replaced.
$ws.Range("E2").Value = "'thinhndp13@uit.edu.vn"
$ws.Range("A2").Value = "1313"
